$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, pushing existing row 11 (and below) down to row 12
$ws.Rows.Item(11).Insert(-4121)  # xlShiftDown = -4121

# Fill in the new row 11 values (vertical JST_GH connector option)
$ws.Cells.Item(11, 1).Value = "JST_GH_V4"
$ws.Cells.Item(11, 2).Value = "BM04B-GHS-TBT"
$ws.Cells.Item(11, 3).Value = "CANV_CAN1, CANV_CAN2"
$ws.Cells.Item(11, 4).Value = "JST_GH_V4"
$ws.Cells.Item(11, 5).Value = "JST_GH_V4"
$ws.Cells.Item(11, 6).Value = ""
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = 2
$ws.Cells.Item(11, 9).Value = "Digi-Key"
$ws.Cells.Item(11, 10).Value = ""
$ws.Cells.Item(11, 12).Value = "455-1580-2-ND"
$ws.Cells.Item(11, 16).Value = ""
